# Update Leve profit-tracking values across all crafting-class sheets.
# Source data refreshed by the scheduled market-price runner (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1125.5892
$ws.Range("J129").Value = 1264.289
$ws.Range("L129").Value = 3792.867
$ws.Range("N129").Value = -13792.867

$ws.Range("H138").Value = 2305.082
$ws.Range("I138").Value = 2956.7856
$ws.Range("J138").Value = 2110.9575
$ws.Range("K138").Value = 8870.356800000001
$ws.Range("L138").Value = 6332.872499999999
$ws.Range("M138").Value = -3730.356800000001
$ws.Range("N138").Value = -16612.8725


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3177.1667
$ws.Range("I61").Value = 2690.1333
$ws.Range("J61").Value = 3988.889
$ws.Range("K61").Value = 2690.1333
$ws.Range("L61").Value = 3988.889
$ws.Range("M61").Value = -2478.1333
$ws.Range("N61").Value = -4412.889

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 30986.97
$ws.Range("I122").Value = 43794.75
$ws.Range("J122").Value = 3042.7273
$ws.Range("K122").Value = 131384.25
$ws.Range("L122").Value = 9128.1819
$ws.Range("M122").Value = -128934.25
$ws.Range("N122").Value = -14028.1819

$ws.Range("H136").Value = 3177.1667
$ws.Range("I136").Value = 2690.1333
$ws.Range("J136").Value = 3988.889
$ws.Range("K136").Value = 8070.3999
$ws.Range("L136").Value = 11966.667
$ws.Range("M136").Value = -5520.3999
$ws.Range("N136").Value = -17066.667


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 53643
$ws.Range("J33").Value = 70024
$ws.Range("L33").Value = 70024
$ws.Range("N33").Value = -70696

$ws.Range("H81").Value = 55361
$ws.Range("J81").Value = 55361
$ws.Range("L81").Value = 55361
$ws.Range("N81").Value = -57483

$ws.Range("H84").Value = 55361
$ws.Range("J84").Value = 55361
$ws.Range("L84").Value = 166083
$ws.Range("N84").Value = -176691

$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws.Range("H132").Value = 42605
$ws.Range("J132").Value = 42605
$ws.Range("L132").Value = 42605
$ws.Range("N132").Value = -52725

$ws.Range("H134").Value = 2858.6086
$ws.Range("I134").Value = 2653.4211
$ws.Range("J134").Value = 3833.25
$ws.Range("K134").Value = 7960.263300000001
$ws.Range("L134").Value = 11499.75
$ws.Range("M134").Value = -5425.263300000001
$ws.Range("N134").Value = -16569.75


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 331.81818
$ws.Range("I22").Value = 225
$ws.Range("J22").Value = 460
$ws.Range("K22").Value = 225
$ws.Range("L22").Value = 460
$ws.Range("M22").Value = 125
$ws.Range("N22").Value = -1160

$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 3000
$ws.Range("K33").Value = 3000
$ws.Range("M33").Value = -2621


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1438.1621
$ws.Range("J39").Value = 1475.3334
$ws.Range("L39").Value = 4426.0002
$ws.Range("N39").Value = -5014.0002

$ws.Range("H69").Value = 22223822
$ws.Range("J69").Value = 22223822
$ws.Range("L69").Value = 66671466
$ws.Range("N69").Value = -66673088

$ws.Range("H72").Value = 22223822
$ws.Range("J72").Value = 22223822
$ws.Range("L72").Value = 200014398
$ws.Range("N72").Value = -200022510

$ws.Range("H82").Value = 2800
$ws.Range("J82").Value = 3166.6667
$ws.Range("L82").Value = 9500.000100000001
$ws.Range("N82").Value = -10312.0001

$ws.Range("H85").Value = 2800
$ws.Range("J85").Value = 3166.6667
$ws.Range("L85").Value = 9500.000100000001
$ws.Range("N85").Value = -12308.0001

$ws.Range("H110").Value = 9968.84
$ws.Range("J110").Value = 10818.762
$ws.Range("L110").Value = 32456.286
$ws.Range("N110").Value = -40636.286

$ws.Range("H113").Value = 1100.4
$ws.Range("I113").Value = 566.6667
$ws.Range("J113").Value = 1194.5883
$ws.Range("K113").Value = 1700.0001
$ws.Range("L113").Value = 3583.7649
$ws.Range("M113").Value = 469.9999
$ws.Range("N113").Value = -7923.7649

$ws.Range("H121").Value = 1058.7441
$ws.Range("I121").Value = 232.8
$ws.Range("J121").Value = 1167.421
$ws.Range("K121").Value = 698.4000000000001
$ws.Range("L121").Value = 3502.263
$ws.Range("M121").Value = 611.5999999999999
$ws.Range("N121").Value = -6122.263

$ws.Range("H134").Value = 7137.375
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 7137.375
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 21412.125
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -31552.125

$ws.Range("H139").Value = 2642.611
$ws.Range("J139").Value = 2778.4285
$ws.Range("L139").Value = 8335.2855
$ws.Range("N139").Value = -18615.2855

$ws.Range("H140").Value = 1406.9048
$ws.Range("I140").Value = 1264.2354
$ws.Range("K140").Value = 3792.7062
$ws.Range("M140").Value = 1387.2938

$ws.Range("H141").Value = 5572.9644
$ws.Range("J141").Value = 7218.5835
$ws.Range("L141").Value = 21655.7505
$ws.Range("N141").Value = -32015.7505


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6044.8184
$ws.Range("I70").Value = 5927
$ws.Range("J70").Value = 6575
$ws.Range("K70").Value = 5927
$ws.Range("L70").Value = 6575
$ws.Range("M70").Value = -5657
$ws.Range("N70").Value = -7115

$ws.Range("H73").Value = 6044.8184
$ws.Range("I73").Value = 5927
$ws.Range("J73").Value = 6575
$ws.Range("K73").Value = 5927
$ws.Range("L73").Value = 6575
$ws.Range("M73").Value = -4991
$ws.Range("N73").Value = -8447

$ws.Range("H122").Value = 3146.5
$ws.Range("I122").Value = 2681.0715
$ws.Range("J122").Value = 3689.5
$ws.Range("K122").Value = 8043.2145
$ws.Range("L122").Value = 11068.5
$ws.Range("M122").Value = -5593.2145
$ws.Range("N122").Value = -15968.5

$ws.Range("H132").Value = 3432.8572
$ws.Range("I132").Value = 3018.5
$ws.Range("J132").Value = 3985.3333
$ws.Range("K132").Value = 9055.5
$ws.Range("L132").Value = 11955.9999
$ws.Range("M132").Value = -6525.5
$ws.Range("N132").Value = -17015.9999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2760.4614
$ws.Range("I122").Value = 2520.6667
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 7562.000100000001
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -5112.000100000001
$ws.Range("N122").Value = -14800

$ws.Range("H132").Value = 3339.4412
$ws.Range("I132").Value = 2488.2273
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 7464.6819
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -4934.6819
$ws.Range("N132").Value = -19760


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1419.4286
$ws.Range("I126").Value = 1290.3636
$ws.Range("J126").Value = 1892.6666
$ws.Range("K126").Value = 3871.0908
$ws.Range("L126").Value = 5677.9998
$ws.Range("M126").Value = -1401.0908
$ws.Range("N126").Value = -10617.9998

$ws.Range("H136").Value = 2660.3438
$ws.Range("I136").Value = 2085.24
$ws.Range("J136").Value = 4714.2856
$ws.Range("K136").Value = 6255.719999999999
$ws.Range("L136").Value = 14142.8568
$ws.Range("M136").Value = -3705.719999999999
$ws.Range("N136").Value = -19242.8568

